$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 changes
$ws.Range("G4").Value = 4.33
$ws.Range("I4").Value = 1.9
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("AI4").Value = 8
$ws.Range("AJ4").Value = 9.5

# Row 5 changes
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62

# Row 13 changes
$ws.Range("Q13").Value = 2.03
$ws.Range("R13").Value = 1.83
